$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force Text format on the Price (D) and Volume(1h) (E) columns for the data rows
# so the numeric-looking / percent-looking strings are preserved verbatim,
# matching the inline-string ("t=inlineStr") cell type used in the source file.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '312.17'
$ws.Range("E2").Value = '0.86%'
$ws.Range("D3").Value = '37.75'
$ws.Range("E3").Value = '0.60%'
$ws.Range("D4").Value = '5.124'
$ws.Range("E4").Value = '0.54%'
$ws.Range("D5").Value = '0.07907'
$ws.Range("E5").Value = '0.69%'
$ws.Range("D6").Value = '4.408'
$ws.Range("E6").Value = '0.82%'
$ws.Range("B7").Value = 'FTXToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D7").Value = '1.902'
$ws.Range("E7").Value = '-3.03%'
$ws.Range("B8").Value = 'KuCoinToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs'
$ws.Range("D8").Value = '8.267'
$ws.Range("E8").Value = '-0.36%'
$ws.Range("D9").Value = '0.9225'
$ws.Range("E9").Value = '-0.45%'
$ws.Range("D10").Value = '0.1206'
$ws.Range("E10").Value = '-10.80%'
$ws.Range("D11").Value = '0.1936'
$ws.Range("E11").Value = '-1.15%'
$ws.Range("D12").Value = '0.09226'
$ws.Range("E12").Value = '3.11%'
$ws.Range("D13").Value = '0.03339'
$ws.Range("E13").Value = '-3.22%'
$ws.Range("D14").Value = '0.09620'
$ws.Range("E14").Value = '-0.86%'
$ws.Range("D15").Value = '0.001386'
$ws.Range("E15").Value = '0.05%'
$ws.Range("D16").Value = '0.005838'
$ws.Range("E16").Value = '-3.17%'
$ws.Range("D17").Value = '3.512'
$ws.Range("E17").Value = '-2.30%'
$ws.Range("D18").Value = '3.099'
$ws.Range("E18").Value = '-1.59%'
$ws.Range("D19").Value = '0.3451'
$ws.Range("E19").Value = '-0.42%'
$ws.Range("D20").Value = '5.284'
$ws.Range("E20").Value = '5.49%'
$ws.Range("E21").Value = '-1.75%'
$ws.Range("E22").Value = '2.89%'
$ws.Range("E23").Value = '-0.20%'
$ws.Range("D24").Value = '0.04361'
$ws.Range("E24").Value = '0.33%'
$ws.Range("E25").Value = '2.41%'
$ws.Range("D26").Value = '0.004312'
$ws.Range("E27").Value = '-9.80%'
$ws.Range("E39").Value = '-6.40%'
$ws.Range("E40").Value = '2.68%'
$ws.Range("D41").Value = '0.007644'
$ws.Range("E41").Value = '0.40%'
$ws.Range("D42").Value = '0.009107'
$ws.Range("E42").Value = '-7.34%'
$ws.Range("E43").Value = '0.58%'
$ws.Range("D44").Value = '0.002010'
$ws.Range("E44").Value = '-2.61%'
$ws.Range("D45").Value = '0.008598'
$ws.Range("D46").Value = '0.00006699'
$ws.Range("E46").Value = '-1.30%'
$ws.Range("E47").Value = '-0.19%'
$ws.Range("D48").Value = '0.001200'
$ws.Range("E48").Value = '-7.80%'
$ws.Range("E49").Value = '-4.23%'
$ws.Range("E50").Value = '-0.19%'
$ws.Range("E51").Value = '-0.19%'
